$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Preceeded" -> "Preceded" in cell D2
$ws.Range("D2").Value = "Preceded by Greenland Stadial 12, ends with GI 11"

# Update the selected cell to D3 (reflecting where the editor's cursor ended up)
$ws.Range("D3").Select()
